$d = $word.ActiveDocument

# Find the "League_Table" paragraph - two new sub-bullets need to be
# inserted immediately before it (right after "League Contact (...)").
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "League_Table`r") {
        $targetIndex = $i
        break
    }
}

$insertionPoint = $d.Paragraphs($targetIndex).Range
$insertionPoint.Collapse(1)

# Insert "Age_group(choice_field)" immediately before "League_Table".
$insertionPoint.InsertParagraphBefore()
$ageRange = $d.Paragraphs($targetIndex).Range
$ageRange.Text = "Age_group(choice_field)"
$ageRange.ListFormat.ListIndent()

# Insert "Gender(choice_field)" immediately after the new Age_group
# paragraph (still before "League_Table").
$targetIndex = $targetIndex + 1
$insertionPoint2 = $d.Paragraphs($targetIndex).Range
$insertionPoint2.Collapse(1)
$insertionPoint2.InsertParagraphBefore()
$genderRange = $d.Paragraphs($targetIndex).Range
$genderRange.Text = "Gender(choice_field)"
$genderRange.ListFormat.ListIndent()

# Relocate the "_GoBack" bookmark so it sits right after the newly
# inserted "Age_group(choice_field)" paragraph instead of after "Points ".
$ageParaIndex = $targetIndex - 1
$ageParaRange = $d.Paragraphs($ageParaIndex).Range
$bmStart = $ageParaRange.End - 1
$bmRange = $d.Range($bmStart, $bmStart)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
